# Add Jorge Jola (and Juan Fiore) to the team roster.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("team")

# Shift existing rows 13-15 (Franco, Roberto, Thatiane) down to 15-17,
# opening up rows 13 and 14 for the two new team members.
$ws.Range("A13:A14").EntireRow.Insert()

# Row 13: Juan Fiore - Visiting Scholar (name + role only, matching target)
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Juan Fiore"
$ws.Cells.Item(13, 3).Value = "Visiting Scholar"

# Row 14: Jorge Jola - Visiting Scholar, with photo, LinkedIn hyperlink, bio
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Jorge Jola"
$ws.Cells.Item(14, 3).Value = "Visiting Scholar"
$ws.Cells.Item(14, 4).Value = "jorge_jola.jpg"
$ws.Cells.Item(14, 5).Value = "https://www.linkedin.com/in/jjola-unal/"
$ws.Cells.Item(14, 6).Value = "Undergraduate student in Agronomic Engineering at the National University of Colombia, with experience in agronomic data analysis. Skilled in statistics, Web applications, and the application of machine learning in agriculture."

# Copy the existing "Hyperlink" cell style (font/underline/color) from another
# linked cell so no new style entries get created in styles.xml.
$ws.Range("E3").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Cells.Item(14, 5).Value = "https://www.linkedin.com/in/jjola-unal/"

# Rebuild hyperlinks: deleting any one hyperlink clears the collection in
# this engine, so remove them all and re-add at the correct final refs.
$ws.Range("E14").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.linkedin.com/in/pedro-cisdeli/")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.linkedin.com/in/ignaciociampitti/")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.linkedin.com/in/leonardo-bosche/")
$ws.Hyperlinks.Add($ws.Range("E9"), "https://www.linkedin.com/in/gmandrini/")
$ws.Hyperlinks.Add($ws.Range("E16"), "https://www.linkedin.com/in/roberto-carlos-romero-palomeque-831917252?utm_source=share&utm_campaign=share_via&utm_content=profile&utm_medium=ios_app")
$ws.Hyperlinks.Add($ws.Range("E14"), "https://www.linkedin.com/in/jjola-unal/", "", "https://www.linkedin.com/in/jjola-unal/", "https://www.linkedin.com/in/jjola-unal/")
